# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Each changed cell is re-written as text (values like prices/percentages are
# stored as strings in this report, not numbers) using a leading apostrophe to
# force text entry, then the cell's style is reset to "Normal" so no stray
# number-format / quote-prefix style is left behind (matching the original,
# un-styled inline-string cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'67.775.99"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.10%  "
$ws.Range('E2').Style = 'Normal'
# Row 3
$ws.Range('D3').Value = "'3.804.87"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.19%  "
$ws.Range('E3').Style = 'Normal'
# Row 4
$ws.Range('E4').Value = "'  +0.15%  "
$ws.Range('E4').Style = 'Normal'
# Row 5
$ws.Range('D5').Value = "'599.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.73%  "
$ws.Range('E5').Style = 'Normal'
# Row 6
$ws.Range('D6').Value = "'167.60"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.57%  "
$ws.Range('E6').Style = 'Normal'
# Row 7
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.04%  "
$ws.Range('E7').Style = 'Normal'
# Row 8
$ws.Range('E8').Value = "'  +0.29%  "
$ws.Range('E8').Style = 'Normal'
# Row 9
$ws.Range('E9').Value = "'  +0.80%  "
$ws.Range('E9').Style = 'Normal'
# Row 10
$ws.Range('E10').Value = "'  -0.94%  "
$ws.Range('E10').Style = 'Normal'
# Row 11
$ws.Range('E11').Value = "'  -0.08%  "
$ws.Range('E11').Style = 'Normal'
# Row 12
$ws.Range('E12').Value = "'  -1.07%  "
$ws.Range('E12').Style = 'Normal'
# Row 13
$ws.Range('D13').Value = "'35.99"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.88%  "
$ws.Range('E13').Style = 'Normal'
# Row 14
$ws.Range('D14').Value = "'4.442.14"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.25%  "
$ws.Range('E14').Style = 'Normal'
# Row 15
$ws.Range('D15').Value = "'3.824.68"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.45%  "
$ws.Range('E15').Style = 'Normal'
# Row 16
$ws.Range('D16').Value = "'18.51"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.02%  "
$ws.Range('E16').Style = 'Normal'
# Row 17
$ws.Range('D17').Value = "'67.853.32"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.23%  "
$ws.Range('E17').Style = 'Normal'
# Row 18
$ws.Range('D18').Value = "'7.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.16%  "
$ws.Range('E18').Style = 'Normal'
# Row 19
$ws.Range('E19').Value = "'  +0.41%  "
$ws.Range('E19').Style = 'Normal'
# Row 20
$ws.Range('D20').Value = "'462.35"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.74%  "
$ws.Range('E20').Style = 'Normal'
# Row 21
$ws.Range('D21').Value = "'9.89"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -3.47%  "
$ws.Range('E21').Style = 'Normal'
# Row 22
$ws.Range('E22').Value = "'  +0.35%  "
$ws.Range('E22').Style = 'Normal'
# Row 23
$ws.Range('D23').Value = "'0.0000151"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -1.49%  "
$ws.Range('E23').Style = 'Normal'
# Row 24
$ws.Range('E24').Value = "'  -0.39%  "
$ws.Range('E24').Style = 'Normal'
# Row 25
$ws.Range('D25').Value = "'12.09"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.43%  "
$ws.Range('E25').Style = 'Normal'
# Row 26
$ws.Range('E26').Value = "'  -1.48%  "
$ws.Range('E26').Style = 'Normal'
# Row 27
$ws.Range('E27').Value = "'  -0.67%  "
$ws.Range('E27').Style = 'Normal'
# Row 28
$ws.Range('E28').Value = "'  -0.10%  "
$ws.Range('E28').Style = 'Normal'
# Row 29
$ws.Range('D29').Value = "'3.953.64"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.20%  "
$ws.Range('E29').Style = 'Normal'
# Row 30
$ws.Range('E30').Value = "'  -0.64%  "
$ws.Range('E30').Style = 'Normal'
# Row 31
$ws.Range('D31').Value = "'7.45"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +2.05%  "
$ws.Range('E31').Style = 'Normal'
# Row 32
$ws.Range('E32').Value = "'  +1.79%  "
$ws.Range('E32').Style = 'Normal'
# Row 33
$ws.Range('D33').Value = "'29.54"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.11%  "
$ws.Range('E33').Style = 'Normal'
# Row 34
$ws.Range('D34').Value = "'1.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.07%  "
$ws.Range('E34').Style = 'Normal'
# Row 36
$ws.Range('D36').Value = "'3.745.37"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.15%  "
$ws.Range('E36').Style = 'Normal'
# Row 37
$ws.Range('D37').Value = "'0.1000"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.35%  "
$ws.Range('E37').Style = 'Normal'
# Row 38
$ws.Range('D38').Value = "'3.42"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +2.06%  "
$ws.Range('E38').Style = 'Normal'
# Row 39
$ws.Range('E39').Value = "'  -0.10%  "
$ws.Range('E39').Style = 'Normal'
# Row 40
$ws.Range('D40').Value = "'1.00"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.18%  "
$ws.Range('E40').Style = 'Normal'
# Row 41
$ws.Range('E41').Value = "'  +0.33%  "
$ws.Range('E41').Style = 'Normal'
# Row 42
$ws.Range('E42').Value = "'  +0.05%  "
$ws.Range('E42').Style = 'Normal'
# Row 43
$ws.Range('E43').Value = "'  -0.01%  "
$ws.Range('E43').Style = 'Normal'
# Row 44
$ws.Range('D44').Value = "'48.12"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +2.17%  "
$ws.Range('E44').Style = 'Normal'
# Row 45
$ws.Range('D45').Value = "'0.301"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.69%  "
$ws.Range('E45').Style = 'Normal'
# Row 46
$ws.Range('B46').Value = "'Arweave"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'42.93"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -4.22%  "
$ws.Range('E46').Style = 'Normal'
# Row 47
$ws.Range('B47').Value = "'EnergySwap"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'27.96"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +10.59%  "
$ws.Range('E47').Style = 'Normal'
# Row 48
$ws.Range('E48').Value = "'  -0.53%  "
$ws.Range('E48').Style = 'Normal'
# Row 49
$ws.Range('E49').Value = "'  +8.75%  "
$ws.Range('E49').Style = 'Normal'
# Row 50
$ws.Range('D50').Value = "'147.95"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.08%  "
$ws.Range('E50').Style = 'Normal'
# Row 51
$ws.Range('E51').Value = "'  +0.30%  "
$ws.Range('E51').Style = 'Normal'
